$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some Price column values look like plain numbers (e.g. "87.30"); the source
# data stores them as text, so force those specific cells to Text format before
# assigning the value to avoid Excel auto-converting them to numeric type.
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'

$ws.Range('D2').Value = '40.007.22'
$ws.Range('D3').Value = '2.208.86'
$ws.Range('E3').Value = '  -0.52%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '294.44'
$ws.Range('E5').Value = '  +1.60%  '
$ws.Range('D6').Value = '87.30'
$ws.Range('E6').Value = '  -0.87%  '
$ws.Range('E7').Value = '  +0.46%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  -0.44%  '
$ws.Range('B10').Value = 'Avalanche'
$ws.Range('C10').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D10').Value = '30.69'
$ws.Range('E10').Value = '  +0.51%  '
$ws.Range('B11').Value = 'OKB'
$ws.Range('C11').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D11').Value = '51.35'
$ws.Range('E11').Value = '  +7.18%  '
$ws.Range('E12').Value = '  +0.14%  '
$ws.Range('E13').Value = '  +2.34%  '
$ws.Range('E14').Value = '  -1.80%  '
$ws.Range('D15').Value = '2.551.75'
$ws.Range('E15').Value = '  -0.52%  '
$ws.Range('D16').Value = '13.79'
$ws.Range('E16').Value = '  -1.38%  '
$ws.Range('D17').Value = '2.168.22'
$ws.Range('E17').Value = '  -2.91%  '
$ws.Range('D18').Value = '0.734'
$ws.Range('E18').Value = '  +0.49%  '
$ws.Range('D19').Value = '39.936.35'
$ws.Range('E19').Value = '  -0.14%  '
$ws.Range('D20').Value = '0.0₃0887'
$ws.Range('E20').Value = '  +0.19%  '
$ws.Range('D21').Value = '11.20'
$ws.Range('E21').Value = '  -3.18%  '
$ws.Range('E22').Value = '  -1.11%  '
$ws.Range('D23').Value = '65.41'
$ws.Range('E23').Value = '  -0.40%  '
$ws.Range('D24').Value = '234.70'
$ws.Range('E24').Value = '  -0.79%  '
$ws.Range('D26').Value = '2.47'
$ws.Range('E26').Value = '  +0.44%  '
$ws.Range('E27').Value = '  -1.31%  '
$ws.Range('D28').Value = '23.09'
$ws.Range('E28').Value = '  +2.23%  '
$ws.Range('E29').Value = '  -4.74%  '
$ws.Range('E30').Value = '  +0.76%  '
$ws.Range('D31').Value = '159.23'
$ws.Range('E31').Value = '  +2.31%  '
$ws.Range('D32').Value = '31.65'
$ws.Range('E32').Value = '  -0.52%  '
$ws.Range('E33').Value = '  -0.04%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').Value = '4.94'
$ws.Range('E34').Value = '  +0.01%  '
$ws.Range('B35').Value = 'LidoDAOToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D35').Value = '3.05'
$ws.Range('E35').Value = '  +6.06%  '
$ws.Range('E36').Value = '  -0.99%  '
$ws.Range('E37').Value = '  -1.24%  '
$ws.Range('E38').Value = '  +1.36%  '
$ws.Range('D39').Value = '0.0999'
$ws.Range('E39').Value = '  +1.26%  '
$ws.Range('E40').Value = '  +2.67%  '
$ws.Range('E41').Value = '  -2.00%  '
$ws.Range('D42').Value = '2.069.50'
$ws.Range('E42').Value = '  -1.54%  '
$ws.Range('E43').Value = '  -2.68%  '
$ws.Range('D44').Value = '19.51'
$ws.Range('E44').Value = '  +10.45%  '
$ws.Range('E45').Value = '  +0.93%  '
$ws.Range('D46').Value = '9.89'
$ws.Range('E46').Value = '  -0.11%  '
$ws.Range('E47').Value = '  +3.28%  '
$ws.Range('E48').Value = '  -9.50%  '
$ws.Range('D49').Value = '2.425.97'
$ws.Range('E49').Value = '  -0.24%  '
$ws.Range('E50').Value = '  +2.11%  '
$ws.Range('E51').Value = '  +0.61%  '

Write-Host "Applied crypto list updates"
